# "updated WA arch with names"
#
# Slide 2 ("CO Device" / "final binary" -> "(FF file)" / "Web App" -> "WA " + "Web App")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1) "Device" -> "CO Device" in the "Device / runtime" rounded-rectangle (id=51) ---
# Shape 4 on this slide holds many blank paragraphs plus "Device" (para 10) and "runtime" (para 11).
$deviceShape = $s.Shapes.Item(4)
$tr = $deviceShape.TextFrame.TextRange
$devicePara = $tr.Paragraphs(10, 1)
# Replace just the word "Device" (first 6 chars of that paragraph) in place so the
# existing run formatting (bold + tx1 solid fill) is preserved and the text becomes
# a single run reading "CO Device".
$devicePara.Characters(1, 6).Text = "CO Device"

# --- 2) Add a new "(FF file)" paragraph under "final binary" (Rectangle 131, id=132) ---
$finalBinaryShape = $s.Shapes.Item(24)
$null = $finalBinaryShape.TextFrame.TextRange.InsertAfter([char]13 + "(FF file)")

# --- 3) Add a new "WA " paragraph above "Web App" (Rectangle 47, id=48) ---
# The shape auto-fits its height (spAutoFit), so inserting the extra line also grows
# the shape from cy=369332 to cy=646331 automatically.
$webAppShape = $s.Shapes.Item(39)
$null = $webAppShape.TextFrame.TextRange.InsertBefore("WA " + [char]13)
